$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate data row (duplicate "OC PATH" / 1492 entry)
$ws.Rows.Item(19).Delete()

# Re-sort the data (A1:B28, including header) descending by column B (Provider Abbrev)
$rng = $ws.Range("A1:B28")
$key = $ws.Range("B1:B28")
$rng.Sort($key, 2, $null, $null, 2, $null, $null, 1)

# Re-apply the AutoFilter over the full table range
$ws.Range("A1:B28").AutoFilter()

# Record the hidden sheet-scoped "_FilterDatabase" defined name (as Excel does for AutoFilter ranges)
$n = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$B`$28")
$n.Visible = $false

# Update the active selection / scroll position to match final view
$ws.Range("B20").Select()
